$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 3852
$ws1.Range("F4").Value = 1378
$ws1.Range("F6").Value = 395
$ws1.Range("F7").Value = 216
$ws1.Range("F8").Value = 64
$ws1.Range("F9").Value = 8873
$ws1.Range("F12").Value = 158
$ws1.Range("F13").Value = 312
$ws1.Range("F14").Value = 354
$ws1.Range("F15").Value = 120
$ws1.Range("F18").Value = 11283
$ws1.Range("F19").Value = 50
$ws1.Range("F20").Value = 299
$ws1.Range("F21").Value = 79
$ws1.Range("F24").Value = 147
$ws1.Range("F29").Value = 2692
$ws1.Range("F35").Value = 924
$ws1.Range("F40").Value = 3037
$ws1.Range("F43").Value = 762
$ws1.Range("F44").Value = 363
$ws1.Range("F45").Value = 367
$ws1.Range("F46").Value = 54
$ws1.Range("F47").Value = 148

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F11").Value = 7
$ws2.Range("F14").Value = 39

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 3852
$ws4.Range("F6").Value = 1378
$ws4.Range("F7").Value = 395
$ws4.Range("F10").Value = 216
$ws4.Range("F11").Value = 64
$ws4.Range("F12").Value = 8873
$ws4.Range("F16").Value = 158
$ws4.Range("F17").Value = 312
$ws4.Range("F18").Value = 354
$ws4.Range("F19").Value = 120
$ws4.Range("F21").Value = 11283
$ws4.Range("F22").Value = 299
$ws4.Range("F25").Value = 147
$ws4.Range("F28").Value = 39
$ws4.Range("F31").Value = 2692
$ws4.Range("F34").Value = 924
$ws4.Range("F39").Value = 3037
$ws4.Range("F43").Value = 363
$ws4.Range("F45").Value = 367
$ws4.Range("F46").Value = 54
$ws4.Range("F47").Value = 148
